$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Effectiveness values (replaces old "Optimistic Scenario" row)
$ws.Range("A2").Value = "Effectiveness "
$ws.Range("B2").Value = "80%"
$ws.Range("C2").Value = "55%"
$ws.Range("D2").Value = "80%"

# Row 3: Duration of Protection values (replaces old "Pessimistic Scenario" row)
$ws.Range("A3").Value = "Duration of Protection"
$ws.Range("B3").Value = "150 days"
$ws.Range("C3").Value = "180 days"
$ws.Range("D3").Value = "2 years "

# Remove old row 4 ("Duration of Protection" row), shifting rows up
$ws.Rows.Item(4).Delete()

# Match the saved selection state from the authored workbook
$ws.Range("J6").Select()
